$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "projectduur ... vermelden." paragraph: drop the gramStart/gramEnd
#    proofErr markers and split the text into two runs ("P" + "rojectduur ...")
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("projectduur wat meer informatie vermelden, precieze data vermelden.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $target = $para.Range
    $xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>P</w:t></w:r><w:r><w:t>rojectduur wat meer informatie vermelden, precieze data vermelden.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $target.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 2) Highlight (green) the "Presentatie: ..." and "Software guidebook: ..."
#    runs.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Presentatie: zoveel tijd van tevoren moet de pp klaar zijn (dit bij meer dingen doen).", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.HighlightColorIndex = 4
}

$rng = $d.Content
$found = $rng.Find.Execute("Software guidebook: bijhouden tijdens het project", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.HighlightColorIndex = 4
}

# ---------------------------------------------------------------------------
# 3) After "... Multilanguage goed aangeven." insert a blank paragraph,
#    "Wireframes in PvA" and "Ja/Nee knoppen".
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Multilanguage goed aangeven.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $nextPara = $para.Next()
    $insertRange = $nextPara.Range
    $insertRange.Collapse(1)
    $xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
<w:p><w:r><w:t xml:space="preserve">Wireframes in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PvA</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
<w:p><w:r><w:t>Ja/Nee knoppen</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $insertRange.InsertXML($xml)
}
